$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11
# from serial date 45206 (2023-10-07) to 45208 (2023-10-09).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
